$d = $word.ActiveDocument

# The first table in the document is the "AdmNo" / "Name" submission table.
# Row 1 = header ("AdmNo"/"Name"), Row 2 = empty, Row 3 = empty (target),
# Row 4 = "2138000"/"Tan Shi Wei Cody", Row 5 = "2112576"/"Li SongLing".
$tbl = $d.Tables.Item(1)

$cellAdm = $tbl.Cell(3, 1)
$cellAdm.Range.Text = "2112688"
$cellAdm.Range.Bold = $true

$cellName = $tbl.Cell(3, 2)
$cellName.Range.Text = "Ho Ka Yee Rachel"
$cellName.Range.Bold = $true
